$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "EndDate" column at H, shifting the existing H:L content
# (ViolationPoint, LocationID, CompanyID, Sertifikasi, Gate) one column to
# the right, to I:M. We copy right-to-left (M<-L, L<-K, K<-J, J<-I, I<-H) so
# we never overwrite a source cell before it has been read, and we use the
# single-argument Copy($dest) form so values, shared-string styles and
# number formats all travel together (reusing existing style indices
# instead of minting new ones).
for ($row = 1; $row -le 3; $row++) {
    $ws.Cells.Item($row, 12).Copy($ws.Cells.Item($row, 13))   # L -> M (Sertifikasi/Gate col)
    $ws.Cells.Item($row, 11).Copy($ws.Cells.Item($row, 12))   # K -> L
    $ws.Cells.Item($row, 10).Copy($ws.Cells.Item($row, 11))   # J -> K
    $ws.Cells.Item($row, 9).Copy($ws.Cells.Item($row, 10))    # I -> J
    $ws.Cells.Item($row, 8).Copy($ws.Cells.Item($row, 9))     # H -> I
}

# New header cell for the inserted column, same style as the other header
# cells in row 1.
$ws.Range("H1").Value = "EndDate"

# New EndDate values for the data rows: reuse the HireDate cell's date
# style/format (column G) and set the same sample date value (45658).
$ws.Cells.Item(2, 7).Copy($ws.Cells.Item(2, 8))
$ws.Range("H2").Value = 45658

$ws.Cells.Item(3, 7).Copy($ws.Cells.Item(3, 8))
$ws.Range("H3").Value = 45658

# Match the row heights shown for the data rows after the edit.
$ws.Rows.Item(2).RowHeight = 34
$ws.Rows.Item(3).RowHeight = 34

# Match the active selection left behind by the edit.
$ws.Range("I7").Select()
